$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark row 7 (Nacho) as Sent ("y")
$ws.Range("C7").Value = "y"

# New row 9: BuffyGirl
$ws.Range("A9").Value = "BuffyGirl"
$ws.Range("B9").Value = "JC Kovacs`n12349 Metric Blvd #1330`nAustin, TX`n78758"
$ws.Range("B9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 60
$ws.Range("C9").Value = "y"

# New row 10: King Darkness
$ws.Range("A10").Value = "King Darkness"
$ws.Range("B10").Value = "Chase Valdez`n2418 East Highway 66`nPMB 539`nGallup NM 87301"
$ws.Range("B10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 60
$ws.Range("C10").Value = "y"

# Update selection to match final state (viewport scrolled to show rows 5+)
$ws.Range("C7").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
